$wb = $excel.ActiveWorkbook

# --- Sheet "Coordenadas" ---
$ws = $wb.Worksheets.Item("Coordenadas")

$ws.Range("C3").Value = 0.004
$ws.Range("D3").Value = -0.004

$ws.Range("C5").Value = 0.005
$ws.Range("D5").Value = -0.004
$ws.Range("E5").Value = 1075.745
$ws.Range("F5").Value = 2103.118

$ws.Range("C7").Value = 0.005
$ws.Range("D7").Value = -0.005
$ws.Range("E7").Value = 1136.195
$ws.Range("F7").Value = 2077.493

# --- Sheet "Parametros Pol" ---
$ws2 = $wb.Worksheets.Item("Parametros Pol")

$ws2.Range("D3").Value = -0.014
$ws2.Range("D4").Value = 0.013
$ws2.Range("D5").Value = 0.0191049731745428
$ws2.Range("D6").Value = 11662.25139205577
